$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- was "Jiovany Ramos", becomes "Aldair Fuentes"
$ws.Range("A2").Value = "Aldair Fuentes"
$ws.Range("B2").Value = "aldair-fuentes"
$ws.Range("C2").Value = "A. Fuentes"
$ws.Range("D2").Value = "M"
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 133
$ws.Range("G2").Value = 876927
$ws.Range("J2").Value = 32.513333333333
$ws.Range("K2").Value = 60.211111111111
$ws.Range("L2").Value = 45

# Row 3 <- was "Aldair Fuentes", becomes "Jiovany Ramos"
$ws.Range("A3").Value = "Jiovany Ramos"
$ws.Range("B3").Value = "ramos-jiovany"
$ws.Range("C3").Value = "J. Ramos"
$ws.Range("D3").Value = "D"
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 1006509
$ws.Range("J3").Value = 40.863636363636
$ws.Range("K3").Value = 20.743636363636
$ws.Range("L3").Value = 55

# Row 4: only F changes
$ws.Range("F4").Value = 93

# Row 6: only F changes
$ws.Range("F6").Value = 473

# Row 7: only F changes
$ws.Range("F7").Value = 324

# Row 8: only F changes
$ws.Range("F8").Value = 83

# Row 9: only F changes
$ws.Range("F9").Value = 271

# Row 10: only F changes
$ws.Range("F10").Value = 98

# Row 11: only F changes
$ws.Range("F11").Value = 329

# Row 12 <- was "Juan Freytes", becomes "Kevin Serna"
$ws.Range("A12").Value = "Kevin Serna"
$ws.Range("B12").Value = "serna-kevin"
$ws.Range("C12").Value = "K. Serna"
$ws.Range("D12").Value = "M"
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 105
$ws.Range("G12").Value = 981374
$ws.Range("J12").Value = 62.135555555556
$ws.Range("K12").Value = 16.653333333333
$ws.Range("L12").Value = 45

# Row 13 <- was "Marco Huaman", becomes "Juan Freytes"
$ws.Range("A13").Value = "Juan Freytes"
$ws.Range("B13").Value = "juan-freytes"
$ws.Range("C13").Value = "J. Freytes"
$ws.Range("D13").Value = "D"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 50
$ws.Range("G13").Value = 962187
$ws.Range("J13").Value = 42.612173913043
$ws.Range("K13").Value = 82.96869565217401
$ws.Range("L13").Value = 115

# Row 14 <- was "Gabriel Costa", becomes "Marco Huaman"
$ws.Range("A14").Value = "Marco Huaman"
$ws.Range("B14").Value = "marco-huaman"
$ws.Range("C14").Value = "M. Huaman"
$ws.Range("D14").Value = "D"
$ws.Range("E14").Value = 25
$ws.Range("F14").Value = 44
$ws.Range("G14").Value = 1090720
$ws.Range("J14").Value = 40.907142857143
$ws.Range("K14").Value = 27.721428571429
$ws.Range("L14").Value = 14

# Row 15 <- was "Kevin Serna", becomes "Gabriel Costa"
$ws.Range("A15").Value = "Gabriel Costa"
$ws.Range("B15").Value = "gabriel-costa"
$ws.Range("C15").Value = "G. Costa"
$ws.Range("D15").Value = "M"
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 585
$ws.Range("G15").Value = 590312
$ws.Range("J15").Value = 61.109302325581
$ws.Range("K15").Value = 52.444186046512
$ws.Range("L15").Value = 43

# Row 16: only F changes
$ws.Range("F16").Value = 86
